# Apply updated cryptocurrency price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.151.86"
$ws.Range("E2").Value = "  -0.94%  "
$ws.Range("D3").Value = "1.861.45"
$ws.Range("E3").Value = "  -0.69%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'0.7088"
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("D6").Value = "'241.35"
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "'0.3108"
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("D9").Value = "'0.07635"
$ws.Range("E9").Value = "  -3.09%  "
$ws.Range("D10").Value = "'24.67"
$ws.Range("E10").Value = "  -2.15%  "
$ws.Range("D11").Value = "'0.08354"
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("D12").Value = "1.856.20"
$ws.Range("E12").Value = "  -1.63%  "
$ws.Range("D13").Value = "'5.191"
$ws.Range("E13").Value = "  -1.89%  "
$ws.Range("D14").Value = "'0.7073"
$ws.Range("E14").Value = "  -3.06%  "
$ws.Range("D15").Value = "'91.08"
$ws.Range("E15").Value = "  -0.20%  "
$ws.Range("D16").Value = "29.170.38"
$ws.Range("D17").Value = "'5.907"
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("D18").Value = "'242.72"
$ws.Range("D19").Value = "'0.000007799"
$ws.Range("E19").Value = "  -0.88%  "
$ws.Range("D20").Value = "2.110.22"
$ws.Range("E20").Value = "  -1.14%  "
$ws.Range("D21").Value = "'13.08"
$ws.Range("E21").Value = "  -2.01%  "
$ws.Range("D22").Value = "'0.9994"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").Value = "'7.861"
$ws.Range("E23").Value = "  -1.52%  "
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("E25").Value = "  -0.84%  "
$ws.Range("D26").Value = "'163.45"
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("D27").Value = "'8.951"
$ws.Range("E27").Value = "  -0.79%  "
$ws.Range("D28").Value = "'18.39"
$ws.Range("E28").Value = "  +0.24%  "
$ws.Range("D29").Value = "'1.322"
$ws.Range("E29").Value = "  -2.95%  "
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").Value = "'4.395"
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("D32").Value = "'4.251"
$ws.Range("E32").Value = "  +3.20%  "
$ws.Range("D33").Value = "'0.05149"
$ws.Range("E33").Value = "  -2.43%  "
$ws.Range("D34").Value = "'0.8008"
$ws.Range("E34").Value = "  +9.84%  "
$ws.Range("D35").Value = "'1.911"
$ws.Range("E35").Value = "  -1.63%  "
$ws.Range("D36").Value = "'1.163"
$ws.Range("E36").Value = "  -2.88%  "
$ws.Range("D37").Value = "'2.687"
$ws.Range("E37").Value = "  +0.52%  "
$ws.Range("D38").Value = "'0.01843"
$ws.Range("E38").Value = "  -1.37%  "
$ws.Range("D39").Value = "'2.692"
$ws.Range("E39").Value = "  -1.15%  "
$ws.Range("D40").Value = "1.164.83"
$ws.Range("E40").Value = "  -5.45%  "
$ws.Range("D41").Value = "'6.207"
$ws.Range("E41").Value = "  +0.27%  "
$ws.Range("D42").Value = "'0.8894"
$ws.Range("E42").Value = "  -2.53%  "
$ws.Range("D43").Value = "'72.81"
$ws.Range("E43").Value = "  -2.75%  "
$ws.Range("D44").Value = "'0.9996"
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").Value = "'102.09"
$ws.Range("E45").Value = "  -0.74%  "
$ws.Range("D46").Value = "2.010.68"
$ws.Range("E46").Value = "  -1.34%  "
$ws.Range("E47").Value = "  -1.33%  "
$ws.Range("D48").Value = "'1.777"
$ws.Range("E48").Value = "  +0.28%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "'0.00000000120"
$ws.Range("E49").Value = "  +2.76%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'9.328"
$ws.Range("E50").Value = "  +0.13%  "
$ws.Range("D51").Value = "'0.4271"
$ws.Range("E51").Value = "  -1.34%  "
